$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 126.291664
$ws.Range("I33").Value = 126.291664
$ws.Range("K33").Value = 126.291664
$ws.Range("M33").Value = 102.708336
$ws.Range("H41").Value = 205.71428
$ws.Range("J41").Value = 474.5
$ws.Range("L41").Value = 474.5
$ws.Range("N41").Value = -1354.5
$ws.Range("H87").Value = 49996.668
$ws.Range("J87").Value = 49996.668
$ws.Range("L87").Value = 49996.668
$ws.Range("N87").Value = -52492.668
$ws.Range("H90").Value = 49996.668
$ws.Range("J90").Value = 49996.668
$ws.Range("L90").Value = 149990.004
$ws.Range("N90").Value = -162470.004
$ws.Range("H113").Value = 8112.6924
$ws.Range("J113").Value = 8308.875
$ws.Range("L113").Value = 8308.875
$ws.Range("N113").Value = -14816.875
$ws.Range("H132").Value = 879.8
$ws.Range("I132").Value = 879.8
$ws.Range("K132").Value = 2639.4
$ws.Range("M132").Value = -109.3999999999996
$ws.Range("H137").Value = 3443.3684
$ws.Range("I137").Value = 2836.9092
$ws.Range("J137").Value = 4277.25
$ws.Range("K137").Value = 8510.7276
$ws.Range("L137").Value = 12831.75
$ws.Range("M137").Value = -5960.7276
$ws.Range("N137").Value = -17931.75
$ws.Range("H138").Value = 2800
$ws.Range("I138").Value = 2333.3333
$ws.Range("K138").Value = 6999.999899999999
$ws.Range("M138").Value = -1859.999899999999
$ws.Range("H141").Value = 2086.889
$ws.Range("I141").Value = 2086.889
$ws.Range("K141").Value = 6260.667
$ws.Range("M141").Value = -1080.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3649.093
$ws.Range("I32").Value = 1422.825
$ws.Range("K32").Value = 1422.825
$ws.Range("M32").Value = -1135.825
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("H88").Value = 3059.9167
$ws.Range("I88").Value = 1005.6667
$ws.Range("J88").Value = 3744.6667
$ws.Range("K88").Value = 1005.6667
$ws.Range("L88").Value = 3744.6667
$ws.Range("M88").Value = -599.6667
$ws.Range("N88").Value = -4556.6667
$ws.Range("H91").Value = 3059.9167
$ws.Range("I91").Value = 1005.6667
$ws.Range("J91").Value = 3744.6667
$ws.Range("K91").Value = 1005.6667
$ws.Range("L91").Value = 3744.6667
$ws.Range("M91").Value = 398.3333
$ws.Range("N91").Value = -6552.6667
$ws.Range("H122").Value = 3042.6667
$ws.Range("I122").Value = 2923
$ws.Range("K122").Value = 8769
$ws.Range("M122").Value = -6319
$ws.Range("H134").Value = 72500
$ws.Range("J134").Value = 72500
$ws.Range("L134").Value = 72500
$ws.Range("N134").Value = -82640
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("H139").Value = 49997
$ws.Range("J139").Value = 49997
$ws.Range("L139").Value = 49997
$ws.Range("N139").Value = -60277

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H134").Value = 11307
$ws.Range("I134").Value = 11281.143
$ws.Range("J134").Value = 11397.5
$ws.Range("K134").Value = 33843.429
$ws.Range("L134").Value = 34192.5
$ws.Range("M134").Value = -31308.429
$ws.Range("N134").Value = -39262.5
$ws.Range("H135").Value = 49998
$ws.Range("J135").Value = 49998
$ws.Range("L135").Value = 49998
$ws.Range("N135").Value = -60138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 13000
$ws.Range("J28").Value = 13000
$ws.Range("L28").Value = 13000
$ws.Range("N28").Value = -13490
$ws.Range("H41").Value = 9822.546
$ws.Range("I41").Value = 1529.5
$ws.Range("J41").Value = 11665.444
$ws.Range("K41").Value = 1529.5
$ws.Range("L41").Value = 11665.444
$ws.Range("M41").Value = -1101.5
$ws.Range("N41").Value = -12521.444
$ws.Range("H43").Value = 10657
$ws.Range("J43").Value = 10657
$ws.Range("L43").Value = 10657
$ws.Range("N43").Value = -11025
$ws.Range("H50").Value = 20180.334
$ws.Range("J50").Value = 19999.8
$ws.Range("L50").Value = 19999.8
$ws.Range("N50").Value = -21249.8
$ws.Range("H101").Value = 10657
$ws.Range("J101").Value = 10657
$ws.Range("L101").Value = 10657
$ws.Range("N101").Value = -17147
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747
$ws.Range("H122").Value = 2147.3333
$ws.Range("I122").Value = 1925.2
$ws.Range("K122").Value = 5775.6
$ws.Range("M122").Value = -3325.6
$ws.Range("H132").Value = 5325
$ws.Range("J132").Value = 5187.5
$ws.Range("L132").Value = 15562.5
$ws.Range("N132").Value = -20622.5
$ws.Range("H134").Value = 4943.5
$ws.Range("I134").Value = 4924.6665
$ws.Range("K134").Value = 14773.9995
$ws.Range("M134").Value = -12238.9995
$ws.Range("H135").Value = 134976
$ws.Range("I135").Value = 70000
$ws.Range("K135").Value = 70000
$ws.Range("M135").Value = -64930

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 900.8333
$ws.Range("I23").Value = 550.5
$ws.Range("J23").Value = 1076
$ws.Range("K23").Value = 1651.5
$ws.Range("L23").Value = 3228
$ws.Range("M23").Value = -1416.5
$ws.Range("N23").Value = -3698
$ws.Range("H34").Value = 1401.25
$ws.Range("J34").Value = 2059.8
$ws.Range("L34").Value = 6179.400000000001
$ws.Range("N34").Value = -6347.400000000001
$ws.Range("H37").Value = 99999.39999999999
$ws.Range("J37").Value = 99999.39999999999
$ws.Range("L37").Value = 299998.2
$ws.Range("N37").Value = -300222.2
$ws.Range("H39").Value = 2499.9167
$ws.Range("J39").Value = 2499.9167
$ws.Range("L39").Value = 7499.750100000001
$ws.Range("N39").Value = -8087.750100000001
$ws.Range("H55").Value = 528.25
$ws.Range("I55").Value = 329.875
$ws.Range("J55").Value = 925
$ws.Range("K55").Value = 989.625
$ws.Range("L55").Value = 2775
$ws.Range("M55").Value = -812.625
$ws.Range("N55").Value = -3129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 16089.083
$ws.Range("I31").Value = 821.125
$ws.Range("J31").Value = 46625
$ws.Range("K31").Value = 821.125
$ws.Range("L31").Value = 46625
$ws.Range("M31").Value = -529.125
$ws.Range("N31").Value = -47209
$ws.Range("H37").Value = 16089.083
$ws.Range("I37").Value = 821.125
$ws.Range("J37").Value = 46625
$ws.Range("K37").Value = 821.125
$ws.Range("L37").Value = 46625
$ws.Range("M37").Value = -544.125
$ws.Range("N37").Value = -47179
$ws.Range("H43").Value = 5080.1763
$ws.Range("I43").Value = 2631.9
$ws.Range("J43").Value = 8577.714
$ws.Range("K43").Value = 2631.9
$ws.Range("L43").Value = 8577.714
$ws.Range("M43").Value = -2480.9
$ws.Range("N43").Value = -8879.714
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("H102").Value = 3247.3076
$ws.Range("I102").Value = 3102.0833
$ws.Range("K102").Value = 3102.0833
$ws.Range("M102").Value = -1480.0833
$ws.Range("H122").Value = 8135.6665
$ws.Range("I122").Value = 10035.923
$ws.Range("J122").Value = 3195
$ws.Range("K122").Value = 30107.769
$ws.Range("L122").Value = 9585
$ws.Range("M122").Value = -27657.769
$ws.Range("N122").Value = -14485
$ws.Range("H132").Value = 6605.25
$ws.Range("I132").Value = 6605.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19815.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17285.75
$ws.Range("N132").Value = ""
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2249.5
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312
$ws.Range("H74").Value = 19598.5
$ws.Range("I74").Value = 19598.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 19598.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -18600.5
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 19598.5
$ws.Range("I77").Value = 19598.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 58795.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -53803.5
$ws.Range("N77").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 28214.84
$ws.Range("J126").Value = 25332.666
$ws.Range("L126").Value = 75997.99800000001
$ws.Range("N126").Value = -80937.99800000001
